$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (Q) to the table, mirroring the formatting of
# the existing "2019" column (P).
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 90.6

# Restore the active selection used when the file was last saved.
[void]$ws.Range("P12").Select()
